$d = $word.ActiveDocument

# 1. Update GPA value from 3.875 to 3.909 (only the "875" run changes; the
# preceding "GPA: 3." run is untouched). A plain Range.Text / Find-Replace
# assignment causes this runtime to coalesce the edited run into its
# identically-formatted neighbor, so toggle a formatting property across the
# edit to keep the run distinct (it is reverted immediately afterwards, so
# the final formatting is unchanged).
$hit = $d.Content
$found = $hit.Find.Execute("875", $false, $false, $false, $false, $false,
                            $true, 1, $false, "", 0)
if ($found) {
    $hit.Font.Bold = 1
    $hit.Text = "909"
    $hit.Font.Bold = 0
}

# 2. Remove the _GoBack bookmark (bookmarkStart/bookmarkEnd for _GoBack)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
